# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500 ...
# Updates the "days-to-process" StructureDefinition summary workbook:
#   - Metadata sheet: URL / Version / Date / Publisher refreshed for the
#     LinuxForHealth rebrand + new IG build.
#   - Elements sheet: the Extension.url row's Fixed Value mirrors the same
#     canonical URL string, so it is refreshed too; and the duplicated
#     ele-1/ext-1 constraint text on the root "Extension" row is cleared
#     (it now only shows once, on Extension.extension).

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/days-to-process"

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = $newUrl
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = $newUrl
$elements.Range("AI2").Value = ""
